# Applies the cryptos.xlsx price/volume/hour refresh described in the commit
# "Updated symbol list on Mon Dec 12 19:17:12 UTC 2022 with GitHub Actions".
# Numeric-looking values (Price column D, Hora column G) are written with a
# leading apostrophe so Excel keeps storing them as text, matching how the
# sheet already represents every other price/volume value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'275.43"
$ws.Range("G2").Value = "'19"

# Row 3
$ws.Range("D3").Value = "'21.08"
$ws.Range("G3").Value = "'19"

# Row 4
$ws.Range("D4").Value = "'6.207"
$ws.Range("G4").Value = "'19"

# Row 5
$ws.Range("D5").Value = "'0.06175"
$ws.Range("G5").Value = "'19"

# Row 6
$ws.Range("D6").Value = "'3.574"
$ws.Range("G6").Value = "'19"

# Row 7
$ws.Range("D7").Value = "'1.511"
$ws.Range("G7").Value = "'19"

# Row 8
$ws.Range("D8").Value = "'6.530"
$ws.Range("G8").Value = "'19"

# Row 9
$ws.Range("D9").Value = "'0.8227"
$ws.Range("G9").Value = "'19"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1646"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("G10").Value = "'19"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08240"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("G11").Value = "'19"

# Row 12
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03418"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G12").Value = "'19"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03146"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("G13").Value = "'19"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09138"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("G14").Value = "'19"

# Row 15
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.773"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("G15").Value = "'19"

# Row 16
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001615"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("G16").Value = "'19"

# Row 17
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04685"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("G17").Value = "'19"

# Row 18
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006450"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("G18").Value = "'19"

# Row 19
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").Value = "'0.006139"
$ws.Range("E19").Value = "18HotbitTokenHTBBestin24h"
$ws.Range("G19").Value = "'19"

# Row 20
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "'0.001068"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("G20").Value = "'19"

# Row 21
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Value = "'0.0001501"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("G21").Value = "'19"

# Row 22
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Value = "'3.725"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("G22").Value = "'19"

# Row 23
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").Value = "'2.318"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("G23").Value = "'19"

# Row 24
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").Value = "'0.01387"
$ws.Range("E24").Value = "23OneONE"
$ws.Range("G24").Value = "'19"

# Row 25
$ws.Range("G25").Value = "'19"

# Row 26
$ws.Range("D26").Value = "'0.1232"
$ws.Range("G26").Value = "'19"

# Row 27
$ws.Range("G27").Value = "'19"

# Row 28
$ws.Range("D28").Value = "'0.0002738"
$ws.Range("G28").Value = "'19"

# Row 29
$ws.Range("G29").Value = "'19"

# Row 30
$ws.Range("G30").Value = "'19"

# Row 31
$ws.Range("G31").Value = "'19"

# Row 32
$ws.Range("G32").Value = "'19"

# Row 33
$ws.Range("G33").Value = "'19"

# Row 34
$ws.Range("G34").Value = "'19"

# Row 35
$ws.Range("G35").Value = "'19"

# Row 36
$ws.Range("G36").Value = "'19"

# Row 37
$ws.Range("G37").Value = "'19"

# Row 38
$ws.Range("G38").Value = "'19"

# Row 39
$ws.Range("G39").Value = "'19"

# Row 40
$ws.Range("D40").Value = "'0.04757"
$ws.Range("G40").Value = "'19"

# Row 41
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.007022"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("G41").Value = "'19"

# Row 42
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1107"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("G42").Value = "'19"

# Row 43
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003521"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("G43").Value = "'19"

# Row 44
$ws.Range("D44").Value = "'0.01044"
$ws.Range("G44").Value = "'19"

# Row 45
$ws.Range("D45").Value = "'0.00005413"
$ws.Range("G45").Value = "'19"

# Row 46
$ws.Range("G46").Value = "'19"

# Row 47
$ws.Range("D47").Value = "'0.7236"
$ws.Range("G47").Value = "'19"

# Row 48
$ws.Range("D48").Value = "'0.001388"
$ws.Range("G48").Value = "'19"

# Row 49
$ws.Range("D49").Value = "'0.00001401"
$ws.Range("G49").Value = "'19"

# Row 50
$ws.Range("D50").Value = "'0.01241"
$ws.Range("G50").Value = "'19"

# Row 51
$ws.Range("G51").Value = "'19"
